$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs / Efna5 / Ephb2 / ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Ephb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.030023666666667
$ws.Range("H2").Value = 6.090071
$ws.Range("I2").Value = 0.8776223887075381
$ws.Range("J2").Value = 0.8776223887075382
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.013267
$ws.Range("N2").Value = 0.039801
$ws.Range("O2").Value = 0.002082344506138891
$ws.Range("P2").Value = 0.002082344506138891
$ws.Range("Q2").Value = 0.02693232398566667
$ws.Range("R2").Value = 0.242390915871
$ws.Range("S2").Value = 0.001827512159589632
$ws.Range("T2").Value = 0.001827512159589632

# Row 3: FAPs / Efna5 / Ephb2 / FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Ephb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.030023666666667
$ws.Range("H3").Value = 6.090071
$ws.Range("I3").Value = 0.8776223887075381
$ws.Range("J3").Value = 0.8776223887075382
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.792929999999999
$ws.Range("N3").Value = 17.37879
$ws.Range("O3").Value = 0.9092391618261221
$ws.Range("P3").Value = 0.9092391618261221
$ws.Range("Q3").Value = 11.75978499934333
$ws.Range("R3").Value = 105.83806499409
$ws.Range("S3").Value = 0.797968645108281
$ws.Range("T3").Value = 0.7979686451082811

# Row 4: FAPs / Efna5 / Ephb2 / M2  (new row)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Ephb2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.030023666666667
$ws.Range("H4").Value = 6.090071
$ws.Range("I4").Value = 0.8776223887075381
$ws.Range("J4").Value = 0.8776223887075382
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05870933333333334
$ws.Range("N4").Value = 0.176128
$ws.Range("O4").Value = 0.009214823074225033
$ws.Range("P4").Value = 0.009214823074225033
$ws.Range("Q4").Value = 0.1191813361208889
$ws.Range("R4").Value = 1.072632025088
$ws.Range("S4").Value = 0.008087135037918712
$ws.Range("T4").Value = 0.008087135037918714

# Row 5: FAPs / Efna5 / Ephb2 / sCs  (new row)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Ephb2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.030023666666667
$ws.Range("H5").Value = 6.090071
$ws.Range("I5").Value = 0.8776223887075381
$ws.Range("J5").Value = 0.8776223887075382
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5062776666666666
$ws.Range("N5").Value = 1.518833
$ws.Range("O5").Value = 0.07946367059351396
$ws.Range("P5").Value = 0.07946367059351396
$ws.Range("Q5").Value = 1.027755645238111
$ws.Range("R5").Value = 9.249800807143
$ws.Range("S5").Value = 0.06973909640174868
$ws.Range("T5").Value = 0.06973909640174868

# Row 6: sCs / Efna5 / Ephb2 / ECs  (new row)
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Ephb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.283071
$ws.Range("H6").Value = 0.849213
$ws.Range("I6").Value = 0.1223776112924619
$ws.Range("J6").Value = 0.1223776112924619
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.013267
$ws.Range("N6").Value = 0.039801
$ws.Range("O6").Value = 0.002082344506138891
$ws.Range("P6").Value = 0.002082344506138891
$ws.Range("Q6").Value = 0.003755502957
$ws.Range("R6").Value = 0.033799526613
$ws.Range("S6").Value = 0.0002548323465492587
$ws.Range("T6").Value = 0.0002548323465492587

# Row 7: sCs / Efna5 / Ephb2 / FAPs  (new row)
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Ephb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.283071
$ws.Range("H7").Value = 0.849213
$ws.Range("I7").Value = 0.1223776112924619
$ws.Range("J7").Value = 0.1223776112924619
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.792929999999999
$ws.Range("N7").Value = 17.37879
$ws.Range("O7").Value = 0.9092391618261221
$ws.Range("P7").Value = 0.9092391618261221
$ws.Range("Q7").Value = 1.63981048803
$ws.Range("R7").Value = 14.75829439227
$ws.Range("S7").Value = 0.111270516717841
$ws.Range("T7").Value = 0.111270516717841

# Row 8: sCs / Efna5 / Ephb2 / M2  (new row)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Ephb2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.283071
$ws.Range("H8").Value = 0.849213
$ws.Range("I8").Value = 0.1223776112924619
$ws.Range("J8").Value = 0.1223776112924619
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05870933333333334
$ws.Range("N8").Value = 0.176128
$ws.Range("O8").Value = 0.009214823074225033
$ws.Range("P8").Value = 0.009214823074225033
$ws.Range("Q8").Value = 0.016618909696
$ws.Range("R8").Value = 0.149570187264
$ws.Range("S8").Value = 0.00112768803630632
$ws.Range("T8").Value = 0.00112768803630632

# Row 9: sCs / Efna5 / Ephb2 / sCs  (new row)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Ephb2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.283071
$ws.Range("H9").Value = 0.849213
$ws.Range("I9").Value = 0.1223776112924619
$ws.Range("J9").Value = 0.1223776112924619
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5062776666666666
$ws.Range("N9").Value = 1.518833
$ws.Range("O9").Value = 0.07946367059351396
$ws.Range("P9").Value = 0.07946367059351396
$ws.Range("Q9").Value = 0.143312525381
$ws.Range("R9").Value = 1.289812728429
$ws.Range("S9").Value = 0.009724574191765284
$ws.Range("T9").Value = 0.009724574191765284
